$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.431.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.686.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.22%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '683.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.686.31'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.20%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.498'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.53%  '
$ws.Range('E10').Value = '  -7.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.27'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.446'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000238'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.306.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.686.84'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.464.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.64'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '472.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.665'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.830.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000129'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -10.57%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.44'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.07%  '
$ws.Range('E30').Value = '  -9.64%  '
$ws.Range('E31').Value = '  -11.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.03'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.166'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.650.74'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.55'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.09'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0938'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.51%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.19%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.955'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '48.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '157.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.91%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -10.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000290'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -10.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '392.40'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -8.07%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '28.38'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.46%  '
